# Fix typo in diagram: the "View for Class N" callouts under the
# generated-Views row on the architecture diagram slide all read
# "View for Class 1" - only the first one is correct; the other three
# should read Class 2, Class 3 and Class 4 respectively.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Shape "Rectangle 27" -> "View for Class 2"
$s.Shapes.Item("Rectangle 27").TextFrame.TextRange.Text = "View for Class 2"

# Shape "Rectangle 28" -> "View for Class 3"
$s.Shapes.Item("Rectangle 28").TextFrame.TextRange.Text = "View for Class 3"

# Shape "Rectangle 29" -> "View for Class 4" (typed as two runs, matching
# how the author appears to have edited just the trailing "1" -> "4")
$tr = $s.Shapes.Item("Rectangle 29").TextFrame.TextRange
$tr.Text = "View for Class 4"
$tr.Characters(10, 7).Font.Size = 8
